# Add bias and activation after convolution
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 ---
# G1, H1, I1 were check-marks ("√"); they become plain numbers now.
$ws.Range("G1").Value = 1
$ws.Range("H1").Value = 3
$ws.Range("I1").Value = 5

# J1, K1 were empty (style 36); now they get text labels and the
# "boxed" style (same visual style already used by I2 / H5).
$ws.Range("I2").Copy() | Out-Null
$ws.Range("J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("I2").Copy() | Out-Null
$ws.Range("K1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("J1").Value = "3+4+5"
$ws.Range("K1").Value = "sum"

# --- Row 2 ---
$ws.Range("G2").Value = 2
$ws.Range("I2").Value = "3+4"

# --- Row 3 ---
$ws.Range("I2").Copy() | Out-Null
$ws.Range("I3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("I3").Value = "1+2+6"

# --- Row 4 ---
$ws.Range("H4").Value = 4

# --- Row 5 ---
$ws.Range("H5").Value = "1+2"

# --- Row 6 ---
# H6 kept referencing the same shared-string slot, but the text itself
# changes from "Final One" to "Final 6".
$ws.Range("H6").Value = "Final 6"

# --- Repaint the fill used by the "boxed" style (now used by I2, H5,
# I3, J1, K1) from black (theme 1) to white (theme 0 / background). ---
foreach ($addr in @("I2","H5","I3","J1","K1")) {
    $ws.Range($addr).Interior.ThemeColor = 2   # msoThemeColorLight1 -> theme="0"
    $ws.Range($addr).Interior.TintAndShade = 0
}

# --- View state: scroll/zoom/selection ---
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("J4").Select() | Out-Null
$excel.ActiveWindow.Zoom = 99
